$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original sheet has rows 2-16 holding one "rotation scheme" per row (column A = numeric
# index, column B = scheme name, columns C:P = averaged intensity values). New rotation
# schemes ("Spiral-...") were run and their results were inserted into the list right after
# "Gaussian-Quadrature", which pushes every later scheme down by three rows. To reproduce
# that row shift (and the corresponding shared-string reordering it causes) we delete the
# affected rows and re-create them from scratch, in the new desired order.

$ws.Range("A10:P16").Delete(-4162) | Out-Null   # xlShiftUp

$schemeNames = @(
  "Gaussian-Quadrature",
  "Spiral-90deg-10rot-5space",
  "Spiral-90deg-15rot-5space",
  "Spiral-90deg-10rot-3space",
  "NoRotation-tilt60deg",
  "Rotation-NoTilt",
  "Rotation-60detTilt",
  "HexGrid-90degTilt5degRes",
  "HexGrid-90degTilt22p5degRes",
  "HexGrid-60degTilt5degRes"
)

$startRow = 10
for ($i = 0; $i -lt $schemeNames.Length; $i++) {
    $row = $startRow + $i
    $idx = 8 + $i

    $cellA = $ws.Cells.Item($row, 1)
    $ws.Range("A2").Copy() | Out-Null
    $cellA.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $cellA.Value = $idx

    $ws.Cells.Item($row, 2).Value = $schemeNames[$i]

    for ($col = 3; $col -le 16; $col++) {
        $ws.Cells.Item($row, $col).Value = 1
    }
}

$excel.CutCopyMode = 0
